# Horarios actualizados Linea 141 - 885
# Re-scrape refresh: new snapshot time 06:46:37 replaces 06:35:33 for most rows
# (a few rows keep their older Hora_Scrap/value because they were not re-touched
# by that particular scrape pass). Each sheet gains extra arrival rows at the end,
# and many of the existing rows below the header shift / get new figures.

$wb = $excel.ActiveWorkbook

# ===================== Sheet "LP1912" =====================
$ws = $wb.Worksheets.Item('LP1912')

# Make room for the 4 new rows scraped at the end of the table
# (shifts the old trailing row(s) down to 75).
$ws.Range("A72:E75").EntireRow.Insert()

# Row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada
$rows = @(
    ,(2, 'Última actualización: 06:46:37', $null, $null, $null, $null)
    ,(3, 'Total filas: 70', $null, $null, $null, $null)
    ,(8, '04:44:46', '04:46', '15_ABASTO', 2, 'LP1912')
    ,(9, '04:44:46', '04:46', '215_EL PELIGRO', 2, 'LP1912')
    ,(10, '03:52:04', '04:46', '215A_EL PATO', 54, 'LP1912')
    ,(44, '06:46:37', '06:46', '17_ROMERO', 0, 'LP1912')
    ,(45, '06:46:37', '06:50', '215A_EL PATO', 4, 'LP1912')
    ,(46, '05:16:02', '06:50', '17_ROMERO', 94, 'LP1912')
    ,(47, '06:35:33', '06:51', '215A_EL PATO', 16, 'LP1912')
    ,(48, '06:46:37', '06:54', '14_ABASTO', 8, 'LP1912')
    ,(49, '06:46:37', '07:04', '225_GOMEZ', 18, 'LP1912')
    ,(50, '06:46:37', '07:06', '215C_EL PATO', 20, 'LP1912')
    ,(51, '06:18:01', '07:07', '215C_EL PATO', 49, 'LP1912')
    ,(52, '06:46:37', '07:13', '14X44_ABASTO', 27, 'LP1912')
    ,(53, '06:18:01', '07:14', '14X44_ABASTO', 56, 'LP1912')
    ,(54, '06:46:37', '07:20', '215A_EL PATO', 34, 'LP1912')
    ,(55, '06:35:33', '07:21', '215A_EL PATO', 46, 'LP1912')
    ,(56, '06:46:37', '07:24', '16_SANTA ANA', 38, 'LP1912')
    ,(57, '06:46:37', '07:29', '14_ABASTO', 43, 'LP1912')
    ,(58, '06:46:37', '07:33', '23_HERNANDEZ', 47, 'LP1912')
    ,(59, '06:46:37', '07:36', '27_EL RETIRO', 50, 'LP1912')
    ,(60, '06:46:37', '07:36', '17X38_ROMERO', 50, 'LP1912')
    ,(61, '06:18:01', '07:37', '27_EL RETIRO', 79, 'LP1912')
    ,(62, '06:46:37', '07:43', '10_OLMOS', 57, 'LP1912')
    ,(63, '06:18:01', '07:44', '10_OLMOS', 86, 'LP1912')
    ,(64, '06:46:37', '07:49', '15_ABASTO', 63, 'LP1912')
    ,(65, '06:35:33', '07:58', '23_HERNANDEZ', 83, 'LP1912')
    ,(66, '06:46:37', '07:59', '11_ETCHEVERRY', 73, 'LP1912')
    ,(67, '06:18:01', '07:59', '23_HERNANDEZ', 101, 'LP1912')
    ,(68, '06:18:01', '08:00', '11_ETCHEVERRY', 102, 'LP1912')
    ,(69, '06:46:37', '08:00', '23_HERNANDEZ', 74, 'LP1912')
    ,(70, '06:46:37', '08:01', '16_SANTA ANA', 75, 'LP1912')
    ,(71, '06:46:37', '08:03', '17X38_ROMERO', 77, 'LP1912')
    ,(72, '06:46:37', '08:14', '10_OLMOS', 88, 'LP1912')
    ,(73, '06:46:37', '08:19', '17_ROMERO', 93, 'LP1912')
    ,(74, '06:46:37', '08:33', '215C_EL PATO', 107, 'LP1912')
    ,(75, '06:35:33', '08:34', '215C_EL PATO', 119, 'LP1912')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 5; $col++) {
        $val = $r[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}

# ===================== Sheet "LP1912-215" =====================
$ws = $wb.Worksheets.Item('LP1912-215')

# Make room for the 2 new rows scraped at the end of the table
# (shifts the old trailing row(s) down to 19).
$ws.Range("A18:E19").EntireRow.Insert()

# Row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada
$rows = @(
    ,(2, 'Última actualización: 06:46:37', $null, $null, $null, $null)
    ,(3, 'Total filas: 14', $null, $null, $null, $null)
    ,(12, '06:46:37', '06:50', '215A_EL PATO', 4, 'LP1912')
    ,(13, '06:35:33', '06:51', '215A_EL PATO', 16, 'LP1912')
    ,(14, '06:46:37', '07:06', '215C_EL PATO', 20, 'LP1912')
    ,(15, '06:18:01', '07:07', '215C_EL PATO', 49, 'LP1912')
    ,(16, '06:46:37', '07:20', '215A_EL PATO', 34, 'LP1912')
    ,(17, '06:35:33', '07:21', '215A_EL PATO', 46, 'LP1912')
    ,(18, '06:46:37', '08:33', '215C_EL PATO', 107, 'LP1912')
    ,(19, '06:35:33', '08:34', '215C_EL PATO', 119, 'LP1912')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 5; $col++) {
        $val = $r[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}

# ===================== Sheet "6203-6173" =====================
$ws = $wb.Worksheets.Item('6203-6173')

# Make room for the 2 new rows scraped at the end of the table
# (shifts the old trailing row(s) down to 11).
$ws.Range("A10:E11").EntireRow.Insert()

# Row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada
$rows = @(
    ,(2, 'Última actualización: 06:46:37', $null, $null, $null, $null)
    ,(3, 'Total filas: 6', $null, $null, $null, $null)
    ,(7, '06:46:37', '07:27', '215A_LA PLATA', 41, 'L6173')
    ,(8, '06:46:37', '08:09', '215A_LA PLATA', 83, 'L6173')
    ,(9, '06:35:33', '08:10', '215A_LA PLATA', 95, 'L6173')
    ,(10, '06:46:37', '08:22', '215C_LA PLATA', 96, 'L6203')
    ,(11, '06:35:33', '08:23', '215C_LA PLATA', 108, 'L6203')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 5; $col++) {
        $val = $r[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}

